$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "33÷2="
$t.Cell(1, 2).Range.Text = "42÷7="
$t.Cell(1, 3).Range.Text = "10÷5="
$t.Cell(1, 4).Range.Text = "89÷4="
$t.Cell(1, 5).Range.Text = "62÷7="
$t.Cell(5, 1).Range.Text = "66÷6="
$t.Cell(5, 2).Range.Text = "92÷9="
$t.Cell(5, 3).Range.Text = "57÷5="
$t.Cell(5, 4).Range.Text = "50÷3="
$t.Cell(5, 5).Range.Text = "15÷8="
$t.Cell(9, 1).Range.Text = "80÷5="
$t.Cell(9, 2).Range.Text = "96÷7="
$t.Cell(9, 3).Range.Text = "82÷5="
$t.Cell(9, 4).Range.Text = "87÷7="
$t.Cell(9, 5).Range.Text = "48÷2="
$t.Cell(13, 1).Range.Text = "11÷2="
$t.Cell(13, 2).Range.Text = "59÷2="
$t.Cell(13, 3).Range.Text = "57÷7="
$t.Cell(13, 4).Range.Text = "92÷9="
$t.Cell(13, 5).Range.Text = "98÷5="
$t.Cell(17, 1).Range.Text = "50÷6="
$t.Cell(17, 2).Range.Text = "60÷2="
$t.Cell(17, 3).Range.Text = "74÷4="
$t.Cell(17, 4).Range.Text = "28÷4="
$t.Cell(17, 5).Range.Text = "56÷2="
